$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("T2").Value = 0.1426048780487805
$ws.Range("V2").Value = 0.0002220611491829204
$ws.Range("Z2").Value = -0.0965673523669755
$ws.Range("AB2").Value = -434.8682906591158
$ws.Range("AC2").Value = "umolO2/min/m2"
$ws.Range("AD2").Value = -434.8682906591158

# Row 3
$ws.Range("T3").Value = 0.1423414634146342
$ws.Range("V3").Value = 0.0002565454225970831
$ws.Range("Z3").Value = -0.1238353147265803
$ws.Range("AB3").Value = -482.7032713075126
$ws.Range("AC3").Value = "umolO2/min/m2"
$ws.Range("AD3").Value = -482.7032713075126

# Row 4
$ws.Range("T4").Value = 0.1477268292682927
$ws.Range("V4").Value = 0.0002491214197856265
$ws.Range("Z4").Value = -0.1496887826809702
$ws.Range("AB4").Value = -600.8667693439612
$ws.Range("AC4").Value = "umolO2/min/m2"
$ws.Range("AD4").Value = -600.8667693439612

# Row 5
$ws.Range("T5").Value = 0.1470341463414634
$ws.Range("V5").Value = 0.0002082674398172554
$ws.Range("Z5").Value = -0.1818591880090642
$ws.Range("AB5").Value = -873.2002859815095
$ws.Range("AC5").Value = "umolO2/min/m2"
$ws.Range("AD5").Value = -873.2002859815095

# Row 6
$ws.Range("T6").Value = 0.1477560975609756
$ws.Range("V6").Value = 0.0001607142857142857
$ws.Range("Z6").Value = -0.1607340688266305
$ws.Range("AB6").Value = -1000.123094921256
$ws.Range("AC6").Value = "umolO2/min/m2"
$ws.Range("AD6").Value = -1000.123094921256

# Row 7
$ws.Range("T7").Value = 0.1455317073170732
$ws.Range("V7").Value = 0.0002247627833421192
$ws.Range("Z7").Value = -0.1285993019301953
$ws.Range("AB7").Value = -572.1556746093943
$ws.Range("AC7").Value = "umolO2/min/m2"
$ws.Range("AD7").Value = -572.1556746093943

# Row 8
$ws.Range("T8").Value = 0.1544
$ws.Range("V8").Value = 0
$ws.Range("Z8").Value = 0.0002688083412225135
$ws.Range("AB8").Value = "Inf"
$ws.Range("AC8").Value = "umolO2/min/m2"
$ws.Range("AD8").Value = "Inf"

# Row 9
$ws.Range("T9").Value = 0.1426048780487805
$ws.Range("V9").Value = 0.0002220611491829204
$ws.Range("Z9").Value = 0.09781887081237563
$ws.Range("AB9").Value = 440.5042087384606
$ws.Range("AC9").Value = "umolO2/min/m2"
$ws.Range("AD9").Value = 440.5042087384606

# Row 10
$ws.Range("T10").Value = 0.1423414634146342
$ws.Range("V10").Value = 0.0002565454225970831
$ws.Range("Z10").Value = 0.09611757129278738
$ws.Range("AB10").Value = 374.6610261830499
$ws.Range("AC10").Value = "umolO2/min/m2"
$ws.Range("AD10").Value = 374.6610261830499

# Row 11
$ws.Range("T11").Value = 0.1477268292682927
$ws.Range("V11").Value = 0.0002491214197856265
$ws.Range("Z11").Value = 0.2250748923543684
$ws.Range("AB11").Value = 903.474669291984
$ws.Range("AC11").Value = "umolO2/min/m2"
$ws.Range("AD11").Value = 903.474669291984

# Row 12
$ws.Range("T12").Value = 0.1470341463414634
$ws.Range("V12").Value = 0.0002082674398172554
$ws.Range("Z12").Value = 0.2005148614679312
$ws.Range("AB12").Value = 962.7758503387439
$ws.Range("AC12").Value = "umolO2/min/m2"
$ws.Range("AD12").Value = 962.7758503387439

# Row 13
$ws.Range("T13").Value = 0.1477560975609756
$ws.Range("V13").Value = 0.0001607142857142857
$ws.Range("Z13").Value = 0.1766729123914516
$ws.Range("AB13").Value = 1099.29812154681
$ws.Range("AC13").Value = "umolO2/min/m2"
$ws.Range("AD13").Value = 1099.29812154681

# Row 14
$ws.Range("T14").Value = 0.1455317073170732
$ws.Range("V14").Value = 0.0002247627833421192
$ws.Range("Z14").Value = 0.05474114188134516
$ws.Range("AB14").Value = 243.550738549192
$ws.Range("AC14").Value = "umolO2/min/m2"
$ws.Range("AD14").Value = 243.550738549192

# Row 15
$ws.Range("T15").Value = 0.1544
$ws.Range("V15").Value = 0
$ws.Range("Z15").Value = 0.001152417643560786
$ws.Range("AB15").Value = "Inf"
$ws.Range("AC15").Value = "umolO2/min/m2"
$ws.Range("AD15").Value = "Inf"
